$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.458.68"
$ws.Range("E2").Value = "  +4.03%  "
$ws.Range("D3").Value = "4.034.12"
$ws.Range("E3").Value = "  +2.91%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "513.97"
$ws.Range("E5").Value = "  -2.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.70"
$ws.Range("E6").Value = "  +1.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.723"
$ws.Range("E7").Value = "  +17.57%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.764"
$ws.Range("E9").Value = "  +5.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.174"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000322"
$ws.Range("E11").Value = "  -4.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.41"
$ws.Range("E12").Value = "  +9.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.77"
$ws.Range("E13").Value = "  +4.52%  "
$ws.Range("D14").Value = "4.666.72"
$ws.Range("E14").Value = "  +2.90%  "
$ws.Range("D15").Value = "4.023.31"
$ws.Range("E15").Value = "  +2.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.07"
$ws.Range("E16").Value = "  +6.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.07"
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.21"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.133"
$ws.Range("E19").Value = "  -1.87%  "
$ws.Range("D20").Value = "72.296.51"
$ws.Range("E20").Value = "  +3.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "442.77"
$ws.Range("E21").Value = "  +2.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "104.11"
$ws.Range("E22").Value = "  +17.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.54"
$ws.Range("E23").Value = "  +4.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.57"
$ws.Range("E24").Value = "  +2.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.96"
$ws.Range("E25").Value = "  -2.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.53"
$ws.Range("E26").Value = "  +1.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.99"
$ws.Range("E27").Value = "  +2.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.62"
$ws.Range("E28").Value = "  +5.94%  "
$ws.Range("E29").Value = "  +2.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.12"
$ws.Range("E30").Value = "  +10.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.66"
$ws.Range("E31").Value = "  +3.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "673.08"
$ws.Range("E32").Value = "  -3.16%  "
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.79"
$ws.Range("E34").Value = "  +14.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "67.08"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "41.24"
$ws.Range("E36").Value = "  +2.70%  "
$ws.Range("B37").Value = "TheGraph"
$ws.Range("C37").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.431"
$ws.Range("E37").Value = "  -2.31%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0861"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("E39").Value = "  +13.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.151"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0488"
$ws.Range("E43").Value = "  +1.27%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.17"
$ws.Range("E44").Value = "  +2.89%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.156"
$ws.Range("E45").Value = "  +11.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.73"
$ws.Range("E46").Value = "  -1.75%  "
$ws.Range("E47").Value = "  +4.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.07"
$ws.Range("E48").Value = "  +1.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.06"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.28"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000267"
$ws.Range("E51").Value = "  +7.35%  "
